$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking values are stored as literal text (matching original inlineStr formatting)
$textCells = @("D2","E2","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","E18","D19","E19","E20","D21","E21","E22","D23","E23","D24","E24","E25","D27","E27","D39","E39","D40","E40","D41","E41","E42","D43","D44","E44","D45","E45","D46","E46","D47","E47","E48","E49","D50","E50","D51","E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply cell value updates
$ws.Range("D2").Value = '308.38'
$ws.Range("E2").Value = '1.08%'
$ws.Range("E3").Value = '1.31%'
$ws.Range("D4").Value = '5.062'
$ws.Range("E4").Value = '0.71%'
$ws.Range("D5").Value = '0.08130'
$ws.Range("E5").Value = '0.60%'
$ws.Range("D6").Value = '2.034'
$ws.Range("E6").Value = '4.31%'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = '4.158'
$ws.Range("E7").Value = '0.50%'
$ws.Range("B8").Value = 'KuCoinToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D8").Value = '7.874'
$ws.Range("E8").Value = '0.36%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '0.9262'
$ws.Range("E9").Value = '-0.52%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '0.1422'
$ws.Range("E10").Value = '13.66%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '0.1935'
$ws.Range("E11").Value = '1.27%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '0.09106'
$ws.Range("E12").Value = '-1.31%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.03455'
$ws.Range("E13").Value = '-1.45%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.09911'
$ws.Range("E14").Value = '-0.36%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '0.001409'
$ws.Range("E15").Value = '-0.86%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '0.006045'
$ws.Range("E16").Value = '-9.70%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '3.836'
$ws.Range("E17").Value = '6.13%'
$ws.Range("E18").Value = '13.07%'
$ws.Range("D19").Value = '0.3459'
$ws.Range("E19").Value = '0.48%'
$ws.Range("E20").Value = '-0.12%'
$ws.Range("D21").Value = '4.810'
$ws.Range("E21").Value = '-6.93%'
$ws.Range("E22").Value = '-7.53%'
$ws.Range("D23").Value = '0.04377'
$ws.Range("E23").Value = '-0.69%'
$ws.Range("D24").Value = '0.001232'
$ws.Range("E24").Value = '-0.27%'
$ws.Range("E25").Value = '4.20%'
$ws.Range("D27").Value = '0.0001299'
$ws.Range("E27").Value = '-0.14%'
$ws.Range("D39").Value = '0.02035'
$ws.Range("E39").Value = '3.72%'
$ws.Range("D40").Value = '0.05150'
$ws.Range("E40").Value = '-0.25%'
$ws.Range("D41").Value = '0.007473'
$ws.Range("E41").Value = '-1.21%'
$ws.Range("E42").Value = '-0.14%'
$ws.Range("D43").Value = '0.1375'
$ws.Range("D44").Value = '0.002128'
$ws.Range("E44").Value = '1.28%'
$ws.Range("D45").Value = '0.009730'
$ws.Range("E45").Value = '-8.63%'
$ws.Range("D46").Value = '0.00006312'
$ws.Range("E46").Value = '-1.05%'
$ws.Range("D47").Value = '0.00000000749'
$ws.Range("E47").Value = '-0.14%'
$ws.Range("E48").Value = '-0.16%'
$ws.Range("E49").Value = '-22.02%'
$ws.Range("D50").Value = '0.00002098'
$ws.Range("E50").Value = '-0.14%'
$ws.Range("D51").Value = '0.0001999'
$ws.Range("E51").Value = '-0.14%'
